$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Insert a new "Jurisdiction" row at position 12 (pushing Description/Purpose/
# Copyright/Immutable down by one), preserving the existing cell formatting (style 2)
# used throughout the table.

# Extend the table formatting down to the new last row (16) by copying row 15's format.
$ws.Range("A15:B15").Copy() | Out-Null
$ws.Range("A16:B16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Shift the old rows 12-15 content down to 13-16 (bottom-up so nothing is clobbered).
for ($r = 15; $r -ge 12; $r--) {
    $a = $ws.Cells.Item($r, 1).Value()
    $b = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r + 1, 1).Value = $a
    $ws.Cells.Item($r + 1, 2).Value = $b
}

# New row 12: Jurisdiction (no value).
$ws.Cells.Item(12, 1).Value = "Jurisdiction"
$ws.Cells.Item(12, 2).Value = ""

# --- Update the metadata values that changed.
$ws.Cells.Item(3, 2).Value = "0.1.7"
$ws.Cells.Item(6, 2).Value = "draft"
$ws.Cells.Item(8, 2).Value = "2024-08-27T12:23:18-05:00"
$ws.Cells.Item(10, 2).Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Cells.Item(11, 2).Value = "Bob Milius (bmilius@nmdp.org)"
